$wb = $excel.ActiveWorkbook

# --- Sheet: Demand_Projection ---
$ws1 = $wb.Worksheets.Item("Demand_Projection")

# row2
$ws1.Range("B2").Value = "ELCCRIXX02"
$ws1.Range("C2").Value = "Output demand of transmission lines in Costa Rica"
$ws1.Range("I2").Value = 40.82
$ws1.Range("J2").Value = 43.12
$ws1.Range("K2").Value = 45.42
$ws1.Range("L2").Value = 47.72
$ws1.Range("M2").Value = 50.02
$ws1.Range("N2").Value = 52.73
$ws1.Range("O2").Value = 55.43
$ws1.Range("P2").Value = 58.14
$ws1.Range("Q2").Value = 60.84
$ws1.Range("R2").Value = 63.55
$ws1.Range("S2").Value = 66.51000000000001
$ws1.Range("T2").Value = 69.48
$ws1.Range("U2").Value = 72.44
$ws1.Range("V2").Value = 75.40000000000001
$ws1.Range("W2").Value = 78.36
$ws1.Range("X2").Value = 81.70999999999999
$ws1.Range("Y2").Value = 85.06
$ws1.Range("Z2").Value = 88.40000000000001
$ws1.Range("AA2").Value = 91.75
$ws1.Range("AB2").Value = 95.09999999999999
$ws1.Range("AC2").Value = 98.75
$ws1.Range("AD2").Value = 102.39
$ws1.Range("AE2").Value = 106.04
$ws1.Range("AF2").Value = 109.68
$ws1.Range("AG2").Value = 113.33
$ws1.Range("AH2").Value = 117.16
$ws1.Range("AI2").Value = 121
$ws1.Range("AJ2").Value = 124.84
$ws1.Range("AK2").Value = 128.67
$ws1.Range("AL2").Value = 132.51

# row3
$ws1.Range("B3").Value = "ELCPANXX02"
$ws1.Range("C3").Value = "Output demand of transmission lines in Panama"
$ws1.Range("I3").Value = 62.21
$ws1.Range("J3").Value = 66.56999999999999
$ws1.Range("K3").Value = 70.93000000000001
$ws1.Range("L3").Value = 75.29000000000001
$ws1.Range("M3").Value = 79.65000000000001
$ws1.Range("N3").Value = 84.16
$ws1.Range("O3").Value = 88.66
$ws1.Range("P3").Value = 93.17
$ws1.Range("Q3").Value = 97.67
$ws1.Range("R3").Value = 102.18
$ws1.Range("S3").Value = 106.88
$ws1.Range("T3").Value = 111.59
$ws1.Range("U3").Value = 116.29
$ws1.Range("V3").Value = 120.99
$ws1.Range("W3").Value = 125.69
$ws1.Range("X3").Value = 130.64
$ws1.Range("Y3").Value = 135.58
$ws1.Range("Z3").Value = 140.52
$ws1.Range("AA3").Value = 145.46
$ws1.Range("AB3").Value = 150.41
$ws1.Range("AC3").Value = 155.74
$ws1.Range("AD3").Value = 161.07
$ws1.Range("AE3").Value = 166.4
$ws1.Range("AF3").Value = 171.73
$ws1.Range("AG3").Value = 177.07
$ws1.Range("AH3").Value = 182.58
$ws1.Range("AI3").Value = 188.1
$ws1.Range("AJ3").Value = 193.62
$ws1.Range("AK3").Value = 199.14
$ws1.Range("AL3").Value = 204.66

# --- Sheet: Profiles ---
$ws2 = $wb.Worksheets.Item("Profiles")

# row2
$ws2.Range("C2").Value = "ELCCRIXX02"
$ws2.Range("D2").Value = "Output demand of transmission lines in Costa Rica"
$ws2.Range("J2").Value = 0.23
$ws2.Range("K2").Value = 0.23
$ws2.Range("L2").Value = 0.23
$ws2.Range("M2").Value = 0.23
$ws2.Range("N2").Value = 0.23
$ws2.Range("O2").Value = 0.23
$ws2.Range("P2").Value = 0.23
$ws2.Range("Q2").Value = 0.23
$ws2.Range("R2").Value = 0.23
$ws2.Range("S2").Value = 0.23
$ws2.Range("T2").Value = 0.23
$ws2.Range("U2").Value = 0.23
$ws2.Range("V2").Value = 0.23
$ws2.Range("W2").Value = 0.23
$ws2.Range("X2").Value = 0.23
$ws2.Range("Y2").Value = 0.23
$ws2.Range("Z2").Value = 0.23
$ws2.Range("AA2").Value = 0.23
$ws2.Range("AB2").Value = 0.23
$ws2.Range("AC2").Value = 0.23
$ws2.Range("AD2").Value = 0.23
$ws2.Range("AE2").Value = 0.23
$ws2.Range("AF2").Value = 0.23
$ws2.Range("AG2").Value = 0.23
$ws2.Range("AH2").Value = 0.23
$ws2.Range("AI2").Value = 0.23
$ws2.Range("AJ2").Value = 0.23
$ws2.Range("AK2").Value = 0.23
$ws2.Range("AL2").Value = 0.23
$ws2.Range("AM2").Value = 0.23

# row3
$ws2.Range("C3").Value = "ELCPANXX02"
$ws2.Range("D3").Value = "Output demand of transmission lines in Panama"
$ws2.Range("J3").Value = 0.21
$ws2.Range("K3").Value = 0.21
$ws2.Range("L3").Value = 0.21
$ws2.Range("M3").Value = 0.21
$ws2.Range("N3").Value = 0.21
$ws2.Range("O3").Value = 0.21
$ws2.Range("P3").Value = 0.21
$ws2.Range("Q3").Value = 0.21
$ws2.Range("R3").Value = 0.21
$ws2.Range("S3").Value = 0.21
$ws2.Range("T3").Value = 0.21
$ws2.Range("U3").Value = 0.21
$ws2.Range("V3").Value = 0.21
$ws2.Range("W3").Value = 0.21
$ws2.Range("X3").Value = 0.21
$ws2.Range("Y3").Value = 0.21
$ws2.Range("Z3").Value = 0.21
$ws2.Range("AA3").Value = 0.21
$ws2.Range("AB3").Value = 0.21
$ws2.Range("AC3").Value = 0.21
$ws2.Range("AD3").Value = 0.21
$ws2.Range("AE3").Value = 0.21
$ws2.Range("AF3").Value = 0.21
$ws2.Range("AG3").Value = 0.21
$ws2.Range("AH3").Value = 0.21
$ws2.Range("AI3").Value = 0.21
$ws2.Range("AJ3").Value = 0.21
$ws2.Range("AK3").Value = 0.21
$ws2.Range("AL3").Value = 0.21
$ws2.Range("AM3").Value = 0.21

# row4
$ws2.Range("C4").Value = "ELCCRIXX02"
$ws2.Range("D4").Value = "Output demand of transmission lines in Costa Rica"
$ws2.Range("J4").Value = 0.27
$ws2.Range("K4").Value = 0.27
$ws2.Range("L4").Value = 0.27
$ws2.Range("M4").Value = 0.27
$ws2.Range("N4").Value = 0.27
$ws2.Range("O4").Value = 0.27
$ws2.Range("P4").Value = 0.27
$ws2.Range("Q4").Value = 0.27
$ws2.Range("R4").Value = 0.27
$ws2.Range("S4").Value = 0.27
$ws2.Range("T4").Value = 0.27
$ws2.Range("U4").Value = 0.27
$ws2.Range("V4").Value = 0.27
$ws2.Range("W4").Value = 0.27
$ws2.Range("X4").Value = 0.27
$ws2.Range("Y4").Value = 0.27
$ws2.Range("Z4").Value = 0.27
$ws2.Range("AA4").Value = 0.27
$ws2.Range("AB4").Value = 0.27
$ws2.Range("AC4").Value = 0.27
$ws2.Range("AD4").Value = 0.27
$ws2.Range("AE4").Value = 0.27
$ws2.Range("AF4").Value = 0.27
$ws2.Range("AG4").Value = 0.27
$ws2.Range("AH4").Value = 0.27
$ws2.Range("AI4").Value = 0.27
$ws2.Range("AJ4").Value = 0.27
$ws2.Range("AK4").Value = 0.27
$ws2.Range("AL4").Value = 0.27
$ws2.Range("AM4").Value = 0.27

# row5
$ws2.Range("C5").Value = "ELCPANXX02"
$ws2.Range("D5").Value = "Output demand of transmission lines in Panama"

# row6
$ws2.Range("C6").Value = "ELCCRIXX02"
$ws2.Range("D6").Value = "Output demand of transmission lines in Costa Rica"
$ws2.Range("J6").Value = 0.23
$ws2.Range("K6").Value = 0.23
$ws2.Range("L6").Value = 0.23
$ws2.Range("M6").Value = 0.23
$ws2.Range("N6").Value = 0.23
$ws2.Range("O6").Value = 0.23
$ws2.Range("P6").Value = 0.23
$ws2.Range("Q6").Value = 0.23
$ws2.Range("R6").Value = 0.23
$ws2.Range("S6").Value = 0.23
$ws2.Range("T6").Value = 0.23
$ws2.Range("U6").Value = 0.23
$ws2.Range("V6").Value = 0.23
$ws2.Range("W6").Value = 0.23
$ws2.Range("X6").Value = 0.23
$ws2.Range("Y6").Value = 0.23
$ws2.Range("Z6").Value = 0.23
$ws2.Range("AA6").Value = 0.23
$ws2.Range("AB6").Value = 0.23
$ws2.Range("AC6").Value = 0.23
$ws2.Range("AD6").Value = 0.23
$ws2.Range("AE6").Value = 0.23
$ws2.Range("AF6").Value = 0.23
$ws2.Range("AG6").Value = 0.23
$ws2.Range("AH6").Value = 0.23
$ws2.Range("AI6").Value = 0.23
$ws2.Range("AJ6").Value = 0.23
$ws2.Range("AK6").Value = 0.23
$ws2.Range("AL6").Value = 0.23
$ws2.Range("AM6").Value = 0.23

# row7
$ws2.Range("C7").Value = "ELCPANXX02"
$ws2.Range("D7").Value = "Output demand of transmission lines in Panama"

# row8
$ws2.Range("C8").Value = "ELCCRIXX02"
$ws2.Range("D8").Value = "Output demand of transmission lines in Costa Rica"
$ws2.Range("J8").Value = 0.27
$ws2.Range("K8").Value = 0.27
$ws2.Range("L8").Value = 0.27
$ws2.Range("M8").Value = 0.27
$ws2.Range("N8").Value = 0.27
$ws2.Range("O8").Value = 0.27
$ws2.Range("P8").Value = 0.27
$ws2.Range("Q8").Value = 0.27
$ws2.Range("R8").Value = 0.27
$ws2.Range("S8").Value = 0.27
$ws2.Range("T8").Value = 0.27
$ws2.Range("U8").Value = 0.27
$ws2.Range("V8").Value = 0.27
$ws2.Range("W8").Value = 0.27
$ws2.Range("X8").Value = 0.27
$ws2.Range("Y8").Value = 0.27
$ws2.Range("Z8").Value = 0.27
$ws2.Range("AA8").Value = 0.27
$ws2.Range("AB8").Value = 0.27
$ws2.Range("AC8").Value = 0.27
$ws2.Range("AD8").Value = 0.27
$ws2.Range("AE8").Value = 0.27
$ws2.Range("AF8").Value = 0.27
$ws2.Range("AG8").Value = 0.27
$ws2.Range("AH8").Value = 0.27
$ws2.Range("AI8").Value = 0.27
$ws2.Range("AJ8").Value = 0.27
$ws2.Range("AK8").Value = 0.27
$ws2.Range("AL8").Value = 0.27
$ws2.Range("AM8").Value = 0.27

# row9
$ws2.Range("C9").Value = "ELCPANXX02"
$ws2.Range("D9").Value = "Output demand of transmission lines in Panama"
$ws2.Range("J9").Value = 0.29
$ws2.Range("K9").Value = 0.29
$ws2.Range("L9").Value = 0.29
$ws2.Range("M9").Value = 0.29
$ws2.Range("N9").Value = 0.29
$ws2.Range("O9").Value = 0.29
$ws2.Range("P9").Value = 0.29
$ws2.Range("Q9").Value = 0.29
$ws2.Range("R9").Value = 0.29
$ws2.Range("S9").Value = 0.29
$ws2.Range("T9").Value = 0.29
$ws2.Range("U9").Value = 0.29
$ws2.Range("V9").Value = 0.29
$ws2.Range("W9").Value = 0.29
$ws2.Range("X9").Value = 0.29
$ws2.Range("Y9").Value = 0.29
$ws2.Range("Z9").Value = 0.29
$ws2.Range("AA9").Value = 0.29
$ws2.Range("AB9").Value = 0.29
$ws2.Range("AC9").Value = 0.29
$ws2.Range("AD9").Value = 0.29
$ws2.Range("AE9").Value = 0.29
$ws2.Range("AF9").Value = 0.29
$ws2.Range("AG9").Value = 0.29
$ws2.Range("AH9").Value = 0.29
$ws2.Range("AI9").Value = 0.29
$ws2.Range("AJ9").Value = 0.29
$ws2.Range("AK9").Value = 0.29
$ws2.Range("AL9").Value = 0.29
$ws2.Range("AM9").Value = 0.29

